# Add a new "Save" column (H) to the s_vals sheet, mirroring the header
# style used by the other header cells in row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1: copy formatting from the existing "sum" header (G1)
# so it picks up the same bold/bordered/centered style, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data column values (H2:H7), plain/unstyled like the other data cells.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 1
